$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Expected" / "Pass" result cells for the QA log + the 3 scenario rows.
# Order matters for shared-string table allocation: "As expected" must be interned
# before "Pass" so they land at indices 33 and 34 respectively.
$ws.Range("F20").Value = "As expected"
$ws.Range("I20").Value = "Pass"

$ws.Range("F21").Value = "As expected"
$ws.Range("I21").Value = "Pass"

$ws.Range("F22").Value = "As expected"
$ws.Range("I22").Value = "Pass"

$ws.Range("J6").Value = "Pass"

# Drop the unused trailing template rows (23-25) from the scenario table.
$ws.Rows("23:25").Delete()

# Restore the selection to match the saved view state.
$ws.Range("N8").Select()
